$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'279.95"
$ws.Range("E2").Value = "'5.24%"
$ws.Range("E3").Value = "'0.87%"
$ws.Range("D4").Value = "'4.950"
$ws.Range("E4").Value = "'5.41%"
$ws.Range("D5").Value = "'0.06408"
$ws.Range("E5").Value = "'5.36%"
$ws.Range("D6").Value = "'7.001"
$ws.Range("E6").Value = "'4.60%"
$ws.Range("D7").Value = "'3.357"
$ws.Range("E7").Value = "'5.89%"
$ws.Range("D8").Value = "'0.8877"
$ws.Range("E8").Value = "'4.43%"
$ws.Range("D9").Value = "'1.042"
$ws.Range("E9").Value = "'15.07%"
$ws.Range("D10").Value = "'0.1495"
$ws.Range("E10").Value = "'6.25%"
$ws.Range("D11").Value = "'0.05171"
$ws.Range("E11").Value = "'4.82%"
$ws.Range("D12").Value = "'0.07359"
$ws.Range("E12").Value = "'3.58%"
$ws.Range("D13").Value = "'0.03149"
$ws.Range("E13").Value = "'0.62%"
$ws.Range("D14").Value = "'0.09078"
$ws.Range("E14").Value = "'0.66%"
$ws.Range("D15").Value = "'0.001551"
$ws.Range("E15").Value = "'1.06%"
$ws.Range("D16").Value = "'0.0006324"
$ws.Range("E16").Value = "'4.07%"
$ws.Range("D17").Value = "'0.006055"
$ws.Range("E17").Value = "'-1.70%"
$ws.Range("D18").Value = "'3.498"
$ws.Range("E18").Value = "'1.15%"
$ws.Range("E19").Value = "'0.83%"
$ws.Range("E20").Value = "'0.79%"
$ws.Range("D21").Value = "'0.1335"
$ws.Range("E21").Value = "'2.67%"
$ws.Range("D22").Value = "'3.930"
$ws.Range("E22").Value = "'-4.16%"
$ws.Range("E23").Value = "'3.06%"
$ws.Range("D24").Value = "'0.001183"
$ws.Range("E24").Value = "'0.30%"
$ws.Range("D25").Value = "'0.003690"
$ws.Range("D26").Value = "'0.0001203"
$ws.Range("E26").Value = "'0.12%"
$ws.Range("D27").Value = "'0.0001700"
$ws.Range("E27").Value = "'1.09%"
$ws.Range("D40").Value = "'0.04094"
$ws.Range("E40").Value = "'4.49%"
$ws.Range("D41").Value = "'0.006666"
$ws.Range("E41").Value = "'59.85%"
$ws.Range("D42").Value = "'0.1181"
$ws.Range("E42").Value = "'6.03%"
$ws.Range("D43").Value = "'0.002365"
$ws.Range("E43").Value = "'12.00%"
$ws.Range("D44").Value = "'0.01255"
$ws.Range("E44").Value = "'8.74%"
$ws.Range("D45").Value = "'0.00005276"
$ws.Range("E45").Value = "'2.86%"
$ws.Range("E46").Value = "'-0.08%"
$ws.Range("E47").Value = "'1,588.11%"
$ws.Range("D48").Value = "'0.02247"
$ws.Range("E48").Value = "'-8.18%"
$ws.Range("E49").Value = "'-0.08%"
$ws.Range("E50").Value = "'-0.15%"
